$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 180 (shifts existing rows 180+ down by one)
$ws.Rows(180).Insert()

# Populate the newly inserted row 180 with the new weekly record
$ws.Cells.Item(180, 1).Value = 11
$ws.Cells.Item(180, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(180, 3).Value = "Bíobío"
$ws.Cells.Item(180, 4).Value = 45001
$ws.Cells.Item(180, 5).Value = 8
$ws.Cells.Item(180, 6).Value = 100112040
$ws.Cells.Item(180, 7).Value = "Cilantro"
$ws.Cells.Item(180, 8).Value = "Sin especificar"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 70
$ws.Cells.Item(180, 11).Value = 7000
$ws.Cells.Item(180, 12).Value = 7500
$ws.Cells.Item(180, 13).Value = 7286
$ws.Cells.Item(180, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(180, 15).Value = "Región Metropolitana"
$ws.Cells.Item(180, 16).Value = 202
$ws.Cells.Item(180, 17).Value = 36
$ws.Cells.Item(180, 18).Value = "Hortaliza"
